$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the current page (bookmark) for "A Student's Guide to Bayesian Statistics"
$ws.Range("C12").Value = 13

# Update the active selection to reflect where the user left off
$ws.Range("E11").Select()
